$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 41

$ws.Range("A$row").Value = "J0VX6D"
$ws.Range("B$row").Value = "Cuchilla de limpieza Ricoh"
$ws.Range("C$row").Value = "Aficio 220 270 1015 1018 1022 1027 1032 2015 2018 2022 2027 2032 3025 3030 3350, MP1800 MP2001 MP2014 MP2352 MP2500 MP2501 MP2510 MP2550 MP2553 MP2851 MP2852 MP3010 MP3053 MP3350 MP3351 MP3352 MP3353"
$ws.Range("D$row").Value = 0
$ws.Range("E$row").Value = 100000
$ws.Range("F$row").Value = 1
$ws.Range("G$row").Value = 0
$ws.Range("H$row").Formula = "=(E$row-D$row)*G$row"
$ws.Range("I$row").Formula = "=D$row*F$row"
$ws.Range("J$row").Value = 0
